# Weekly update: insert a new week's worth of data (2 rows: Primera / Segunda)
# at the top of the "Locoto" price history table, pushing the existing rows
# (37-68) down by two to become rows 39-70.
#
# The new rows re-use the same fixed metadata (Mercado, Region, Categoria,
# Variedad, Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion)
# as the rows that used to sit at 37/38 (now at 39/40 after the shift), and
# only the date / volume / price columns change for the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 37:68 down to 39:70, opening up two blank rows.
$ws.Rows("37:38").Insert()

# Seed the two new rows with a copy of the (now shifted) rows right below them
# so all the unchanged columns (Mercado, Region, Categoria, Variedad, Calidad,
# Unidad, Origen, Kg o Unidades, Clasificacion, etc.) carry over correctly.
$ws.Rows(39).Copy()
$ws.Rows(37).PasteSpecial()

$ws.Rows(40).Copy()
$ws.Rows(38).PasteSpecial()

# Row 37: Locoto, Primera, week of 2021-09-13
$ws.Cells.Item(37, 4).Value = 44452
$ws.Cells.Item(37, 10).Value = 120
$ws.Cells.Item(37, 11).Value = 25000
$ws.Cells.Item(37, 12).Value = 26000
$ws.Cells.Item(37, 13).Value = 25500
$ws.Cells.Item(37, 16).Value = 1275

# Row 38: Locoto, Segunda, week of 2021-09-13
$ws.Cells.Item(38, 4).Value = 44452
$ws.Cells.Item(38, 10).Value = 120
$ws.Cells.Item(38, 11).Value = 22000
$ws.Cells.Item(38, 12).Value = 23000
$ws.Cells.Item(38, 13).Value = 22500
$ws.Cells.Item(38, 16).Value = 1125
